$wb = $excel.ActiveWorkbook

# Update the "展览" sheet (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4960
$ws1.Range("F4").Value = 869

# Update the "全部类型" sheet (all types), which mirrors the same data
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4960
$ws4.Range("F4").Value = 869
